$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated results for Case_2_122 (380 kV case) - vm_pu values
# Each row corresponds to a bus (row r=2..25 -> bus index 0..23)

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.021075889624625
$ws.Cells.Item(2, 4).Value = 1.025362066925665
$ws.Cells.Item(2, 5).Value = 1.046575217321456
$ws.Cells.Item(2, 6).Value = 1.050310168808472
$ws.Cells.Item(2, 9).Value = 1.027980029273226
$ws.Cells.Item(2, 10).Value = 1.026269399659629
$ws.Cells.Item(2, 11).Value = 1.028187860595328
$ws.Cells.Item(2, 12).Value = 1.049340269291631
$ws.Cells.Item(2, 13).Value = 1.053064801641427
$ws.Cells.Item(2, 14).Value = 1.027726819921891

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.022031439966696
$ws.Cells.Item(3, 4).Value = 1.026040804710364
$ws.Cells.Item(3, 5).Value = 1.047855161132216
$ws.Cells.Item(3, 6).Value = 1.051667638027931
$ws.Cells.Item(3, 9).Value = 1.028124167480241
$ws.Cells.Item(3, 10).Value = 1.026862381116127
$ws.Cells.Item(3, 11).Value = 1.028674234398509
$ws.Cells.Item(3, 12).Value = 1.050430667270378
$ws.Cells.Item(3, 13).Value = 1.054233287173958
$ws.Cells.Item(3, 14).Value = 1.028320643480074

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.022649248253463
$ws.Cells.Item(4, 4).Value = 1.026479088988873
$ws.Cells.Item(4, 5).Value = 1.04868391329289
$ws.Cells.Item(4, 6).Value = 1.052546475299208
$ws.Cells.Item(4, 9).Value = 1.028215222617617
$ws.Cells.Item(4, 10).Value = 1.027245005041265
$ws.Cells.Item(4, 11).Value = 1.028987371963658
$ws.Cells.Item(4, 12).Value = 1.051136219262722
$ws.Cells.Item(4, 13).Value = 1.054989303450709
$ws.Cells.Item(4, 14).Value = 1.028703810775073

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.022908855864591
$ws.Cells.Item(5, 4).Value = 1.026663126212781
$ws.Cells.Item(5, 5).Value = 1.049032452619543
$ws.Cells.Item(5, 6).Value = 1.05291605156019
$ws.Cells.Item(5, 9).Value = 1.028252972082331
$ws.Cells.Item(5, 10).Value = 1.027405602504602
$ws.Cells.Item(5, 11).Value = 1.029118636363912
$ws.Cells.Item(5, 12).Value = 1.051432832982444
$ws.Cells.Item(5, 13).Value = 1.055307117109532
$ws.Cells.Item(5, 14).Value = 1.028864636305229

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.022952438144353
$ws.Cells.Item(6, 4).Value = 1.026694014090774
$ws.Cells.Item(6, 5).Value = 1.049090981698313
$ws.Cells.Item(6, 6).Value = 1.052978111729652
$ws.Cells.Item(6, 9).Value = 1.028259279280386
$ws.Cells.Item(6, 10).Value = 1.027432552407368
$ws.Cells.Item(6, 11).Value = 1.029140654011613
$ws.Cells.Item(6, 12).Value = 1.051482635760071
$ws.Cells.Item(6, 13).Value = 1.055360478555462
$ws.Cells.Item(6, 14).Value = 1.028891624479948

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.022652717610761
$ws.Cells.Item(7, 4).Value = 1.026481548959085
$ws.Cells.Item(7, 5).Value = 1.048688569971656
$ws.Cells.Item(7, 6).Value = 1.052551413148492
$ws.Cells.Item(7, 9).Value = 1.028215729111551
$ws.Cells.Item(7, 10).Value = 1.027247151965218
$ws.Cells.Item(7, 11).Value = 1.028989127413632
$ws.Cells.Item(7, 12).Value = 1.051140182628998
$ws.Cells.Item(7, 13).Value = 1.054993550151551
$ws.Cells.Item(7, 14).Value = 1.028705960747904

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.021398925129132
$ws.Cells.Item(8, 4).Value = 1.025591636239618
$ws.Cells.Item(8, 5).Value = 1.047007669076055
$ws.Cells.Item(8, 6).Value = 1.050768837040347
$ws.Cells.Item(8, 9).Value = 1.028029199104545
$ws.Cells.Item(8, 10).Value = 1.026470022846452
$ws.Cells.Item(8, 11).Value = 1.028352559408285
$ws.Cells.Item(8, 12).Value = 1.049708777144681
$ws.Cells.Item(8, 13).Value = 1.053459712700611
$ws.Cells.Item(8, 14).Value = 1.027927728016651

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.019185781291819
$ws.Cells.Item(9, 4).Value = 1.024016610580602
$ws.Cells.Item(9, 5).Value = 1.044049778065955
$ws.Cells.Item(9, 6).Value = 1.047631172954742
$ws.Cells.Item(9, 9).Value = 1.027683593521821
$ws.Cells.Item(9, 10).Value = 1.025092405115136
$ws.Cells.Item(9, 11).Value = 1.027218781733771
$ws.Cells.Item(9, 12).Value = 1.047186318117073
$ws.Cells.Item(9, 13).Value = 1.050756269410039
$ws.Cells.Item(9, 14).Value = 1.026548153910142

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.017707795927658
$ws.Cells.Item(10, 4).Value = 1.02296201691735
$ws.Cells.Item(10, 5).Value = 1.042080464121078
$ws.Cells.Item(10, 6).Value = 1.045541596457992
$ws.Cells.Item(10, 9).Value = 1.027441848423023
$ws.Cells.Item(10, 10).Value = 1.024168485045393
$ws.Cells.Item(10, 11).Value = 1.02645485471853
$ws.Cells.Item(10, 12).Value = 1.045504477662009
$ws.Cells.Item(10, 13).Value = 1.048953446877787
$ws.Cells.Item(10, 14).Value = 1.025622921767925

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.017067204313006
$ws.Cells.Item(11, 4).Value = 1.02250429107664
$ws.Cells.Item(11, 5).Value = 1.04122832109892
$ws.Cells.Item(11, 6).Value = 1.044637279102023
$ws.Cells.Item(11, 9).Value = 1.027334486340885
$ws.Cells.Item(11, 10).Value = 1.023767113526774
$ws.Cells.Item(11, 11).Value = 1.026122157257685
$ws.Cells.Item(11, 12).Value = 1.044776150001431
$ws.Cells.Item(11, 13).Value = 1.048172655415035
$ws.Cells.Item(11, 14).Value = 1.025220980255711

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.016829167532578
$ws.Cells.Item(12, 4).Value = 1.0223341098291
$ws.Cells.Item(12, 5).Value = 1.040911882302744
$ws.Cells.Item(12, 6).Value = 1.044301445321868
$ws.Cells.Item(12, 9).Value = 1.027294204344689
$ws.Cells.Item(12, 10).Value = 1.023617829855535
$ws.Cells.Item(12, 11).Value = 1.02599829178995
$ws.Cells.Item(12, 12).Value = 1.04450560252284
$ws.Cells.Item(12, 13).Value = 1.047882608854049
$ws.Cells.Item(12, 14).Value = 1.025071484584535

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.016880231421524
$ws.Cells.Item(13, 4).Value = 1.022370621580419
$ws.Cells.Item(13, 5).Value = 1.040979755689764
$ws.Cells.Item(13, 6).Value = 1.044373479684499
$ws.Cells.Item(13, 9).Value = 1.027302863196618
$ws.Cells.Item(13, 10).Value = 1.023649860596635
$ws.Cells.Item(13, 11).Value = 1.026024874325386
$ws.Cells.Item(13, 12).Value = 1.044563636551234
$ws.Cells.Item(13, 13).Value = 1.047944826030949
$ws.Cells.Item(13, 14).Value = 1.025103560812963

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.017047530007631
$ws.Cells.Item(14, 4).Value = 1.022490227131749
$ws.Cells.Item(14, 5).Value = 1.041202162448581
$ws.Cells.Item(14, 6).Value = 1.044609517565214
$ws.Cells.Item(14, 9).Value = 1.027331164835882
$ws.Cells.Item(14, 10).Value = 1.02375477769828
$ws.Cells.Item(14, 11).Value = 1.026111924347963
$ws.Cells.Item(14, 12).Value = 1.044753786763931
$ws.Cells.Item(14, 13).Value = 1.048148680615827
$ws.Cells.Item(14, 14).Value = 1.025208626908926

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.017150595918803
$ws.Cells.Item(15, 4).Value = 1.022563898682307
$ws.Cells.Item(15, 5).Value = 1.041339205789619
$ws.Cells.Item(15, 6).Value = 1.04475495748432
$ws.Cells.Item(15, 9).Value = 1.027348549029351
$ws.Cells.Item(15, 10).Value = 1.023819394562758
$ws.Cells.Item(15, 11).Value = 1.026165520746784
$ws.Cells.Item(15, 12).Value = 1.044870942644225
$ws.Cells.Item(15, 13).Value = 1.048274278662804
$ws.Cells.Item(15, 14).Value = 1.025273335536762

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.017750296921316
$ws.Cells.Item(16, 4).Value = 1.022992371986897
$ws.Cells.Item(16, 5).Value = 1.04213703011396
$ws.Cells.Item(16, 6).Value = 1.045601622976282
$ws.Cells.Item(16, 9).Value = 1.027448917160966
$ws.Cells.Item(16, 10).Value = 1.024195095201164
$ws.Cells.Item(16, 11).Value = 1.026476894507256
$ws.Cells.Item(16, 12).Value = 1.045552812505205
$ws.Cells.Item(16, 13).Value = 1.049005261931772
$ws.Cells.Item(16, 14).Value = 1.025649569713169

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.018126308900853
$ws.Cells.Item(17, 4).Value = 1.023260853399593
$ws.Cells.Item(17, 5).Value = 1.042637638525161
$ws.Cells.Item(17, 6).Value = 1.04613284134947
$ws.Cells.Item(17, 9).Value = 1.027511156932928
$ws.Cells.Item(17, 10).Value = 1.024430412162059
$ws.Cells.Item(17, 11).Value = 1.026671699455037
$ws.Cells.Item(17, 12).Value = 1.045980508586213
$ws.Cells.Item(17, 13).Value = 1.049463744731829
$ws.Cells.Item(17, 14).Value = 1.02588522085114

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.018345571083418
$ws.Cells.Item(18, 4).Value = 1.023417349828975
$ws.Cells.Item(18, 5).Value = 1.042929691353049
$ws.Cells.Item(18, 6).Value = 1.046442738823255
$ws.Cells.Item(18, 9).Value = 1.027547201260591
$ws.Cells.Item(18, 10).Value = 1.024567542352032
$ws.Cells.Item(18, 11).Value = 1.026785141373951
$ws.Cells.Item(18, 12).Value = 1.046229968936803
$ws.Cells.Item(18, 13).Value = 1.049731154952189
$ws.Cells.Item(18, 14).Value = 1.026022545781712

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.018420323780328
$ws.Cells.Item(19, 4).Value = 1.023470693364078
$ws.Cells.Item(19, 5).Value = 1.043029283464703
$ws.Cells.Item(19, 6).Value = 1.046548413868462
$ws.Cells.Item(19, 9).Value = 1.02755944748741
$ws.Cells.Item(19, 10).Value = 1.024614278792123
$ws.Cells.Item(19, 11).Value = 1.026823790826721
$ws.Cells.Item(19, 12).Value = 1.046315027257829
$ws.Cells.Item(19, 13).Value = 1.049822332442568
$ws.Cells.Item(19, 14).Value = 1.026069348592908

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.01808597250348
$ws.Cells.Item(20, 4).Value = 1.023232058661052
$ws.Cells.Item(20, 5).Value = 1.0425839221633
$ws.Cells.Item(20, 6).Value = 1.046075841824926
$ws.Cells.Item(20, 9).Value = 1.027504505988919
$ws.Cells.Item(20, 10).Value = 1.024405177933024
$ws.Cells.Item(20, 11).Value = 1.026650817809298
$ws.Cells.Item(20, 12).Value = 1.04593462163957
$ws.Cells.Item(20, 13).Value = 1.049414555437653
$ws.Cells.Item(20, 14).Value = 1.025859950786605

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.016998267300286
$ws.Cells.Item(21, 4).Value = 1.022455010730013
$ws.Cells.Item(21, 5).Value = 1.041136666873468
$ws.Cells.Item(21, 6).Value = 1.044540008382551
$ws.Cells.Item(21, 9).Value = 1.027322841828184
$ws.Cells.Item(21, 10).Value = 1.023723887648199
$ws.Cells.Item(21, 11).Value = 1.02608629820006
$ws.Cells.Item(21, 12).Value = 1.04469779267407
$ws.Cells.Item(21, 13).Value = 1.048088651252705
$ws.Cells.Item(21, 14).Value = 1.025177692991431

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.016313848369878
$ws.Cells.Item(22, 4).Value = 1.021965515878661
$ws.Cells.Item(22, 5).Value = 1.040227210938034
$ws.Cells.Item(22, 6).Value = 1.043574772349739
$ws.Cells.Item(22, 9).Value = 1.02720629128776
$ws.Cells.Item(22, 10).Value = 1.023294397031885
$ws.Cells.Item(22, 11).Value = 1.025729703087814
$ws.Cells.Item(22, 12).Value = 1.043920066220185
$ws.Cells.Item(22, 13).Value = 1.047254852031779
$ws.Cells.Item(22, 14).Value = 1.024747592449178

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.016676722641924
$ws.Cells.Item(23, 4).Value = 1.022225094647734
$ws.Cells.Item(23, 5).Value = 1.040709285044842
$ws.Cells.Item(23, 6).Value = 1.044086424872822
$ws.Cells.Item(23, 9).Value = 1.027268297751917
$ws.Cells.Item(23, 10).Value = 1.02352218581892
$ws.Cells.Item(23, 11).Value = 1.025918898035498
$ws.Cells.Item(23, 12).Value = 1.044332362328167
$ws.Cells.Item(23, 13).Value = 1.047696879788893
$ws.Cells.Item(23, 14).Value = 1.024975704722419

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.018104198968022
$ws.Cells.Item(24, 4).Value = 1.023245070085985
$ws.Cells.Item(24, 5).Value = 1.042608194100086
$ws.Cells.Item(24, 6).Value = 1.046101597311309
$ws.Cells.Item(24, 9).Value = 1.027507512064716
$ws.Cells.Item(24, 10).Value = 1.024416580584777
$ws.Cells.Item(24, 11).Value = 1.026660253896207
$ws.Cells.Item(24, 12).Value = 1.045955355998292
$ws.Cells.Item(24, 13).Value = 1.049436782008236
$ws.Cells.Item(24, 14).Value = 1.025871369631432

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.019758382488281
$ws.Cells.Item(25, 4).Value = 1.024424602935537
$ws.Cells.Item(25, 5).Value = 1.044813994078941
$ws.Cells.Item(25, 6).Value = 1.048441937200393
$ws.Cells.Item(25, 9).Value = 1.027774942139128
$ws.Cells.Item(25, 10).Value = 1.025449524338988
$ws.Cells.Item(25, 11).Value = 1.027513316304579
$ws.Cells.Item(25, 12).Value = 1.0478384618179
$ws.Cells.Item(25, 13).Value = 1.051455259382623
$ws.Cells.Item(25, 14).Value = 1.026905780284254

Write-Output "vm_pu values updated for 380 kV case"